$wb = $excel.ActiveWorkbook

# "train" sheet: bold the header row (A1:C1)
$wsTrain = $wb.Worksheets.Item("train")
$wsTrain.Range("A1:C1").Font.Bold = $true

# "config" sheet: bold the header row (A1:B1) and fix the train_iteration value
$wsConfig = $wb.Worksheets.Item("config")
$wsConfig.Range("A1:B1").Font.Bold = $true
$wsConfig.Range("B5").Value = 20
